$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header row (row 1) so headers use the "<formatversion>"
#    suffix instead of the old "_old" / "_new" suffix:
#      *_old -> *_FV2404
#      *_new -> *_FV2410
# ---------------------------------------------------------------------------
$headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404",
    "diff",
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# ---------------------------------------------------------------------------
# 2) Turn A1:U80 into an Excel Table (ListObject named "Table1").
#
#    The header row (A1:U1) already carries explicit formatting (bold font,
#    grey fill, border). If a ListObject is created directly on top of it,
#    the engine captures that pre-existing look as a table "header row" dxf
#    override (headerRowDxfId + a new <dxf> entry in styles.xml) - which the
#    target workbook does not have (dxfs stays at count="0").
#
#    To avoid that, the table is first created on a throwaway, unformatted
#    range, then resized onto the real A1:U80 range - resizing does not
#    re-capture formatting. The throwaway header values are cleaned up
#    afterwards and the real header text is (re)written directly into
#    A1:U1, which also updates the table's column names.
# ---------------------------------------------------------------------------
$tmpRow1 = 200
$tmpRow2 = 201
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item($tmpRow1, $c).Value = "TmpCol" + $c
}
$ws.Cells.Item($tmpRow2, 1).Value = "x"

$tmpRange = $ws.Range($ws.Cells.Item($tmpRow1, 1), $ws.Cells.Item($tmpRow2, 21))
$tbl = $ws.ListObjects.Add(1, $tmpRange, $null, 1)
$tbl.Name = "Table1"

$realRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item(80, 21))
$tbl.Resize($realRange)

# Remove the throwaway rows now that the table lives on the real range.
$ws.Range($ws.Cells.Item($tmpRow1, 1), $ws.Cells.Item($tmpRow2, 21)).Clear()

# Write the real header text - this also renames the table's columns.
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# No explicit table style (matches the target's bare <tableStyleInfo/>).
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowAutoFilter = $true

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1, frozen).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
